$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match the date cell formatting already used by row 3 (A3) for the new rows
$ws.Range("A3").Copy()
$ws.Range("A4:A5").PasteSpecial(-4122)  # xlPasteFormats

# Row 4: date 2017-01-31, hours 1, activity "ER Diagram"
$ws.Range("A4").Value2 = 42766
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = "ER Diagram"

# Row 5: date 2017-02-01, hours 1, activity "ER Diagram"
$ws.Range("A5").Value2 = 42767
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "ER Diagram"

# Update the active selection to C6
$ws.Range("C6").Select()
